$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.238.93"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.589.22"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'212.00"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.245"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").Value = "'19.22"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.811.87"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.616.29"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "'63.84"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "26.236.73"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'7.43"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'214.01"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").Value = "'144.52"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'15.12"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'0.0495"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "1.420.18"
$ws.Range("E33").Value = "  +8.56%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'0.590"
$ws.Range("E36").Value = "  -4.06%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "'5.91"
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("D40").Value = "'0.822"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'0.938"
$ws.Range("E42").Value = "  -14.44%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.723.93"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "'61.17"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.49"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0502"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0966"
$ws.Range("E51").Value = "  -1.55%  "
